# Add a "pagecount" column (I) with per-book page counts, matching the
# commit "feat: display page count and genres".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column I.
$ws.Range("I1").Value = "pagecount"

# Page counts for rows 2..38 (book id 1..37), in row order.
$pagecounts = @(319,297,175,215,417,434,310,184,416,172,255,265,244,238,356,328,247,255,353,281,425,306,330,417,226,421,340,136,297,614,427,580,306,354,370,255,449)

for ($i = 0; $i -lt $pagecounts.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $pagecounts[$i]
}

# Give the new column (and its neighbours G/H, which previously relied on
# the default width) explicit widths, matching the new layout as closely as
# Excel's character-based column-width model allows.
$ws.Columns.Item(7).ColumnWidth = 8.6    # G: type
$ws.Columns.Item(8).ColumnWidth = 25.3   # H: genres
$ws.Columns.Item(9).ColumnWidth = 8.6    # I: pagecount

# Restore the author's scroll position / selection on the sheet.
$ws.Range("D35").Select()
